$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the account rows that were dropped from the export:
#   row 2  -> 004213929 RODOLFO   250083.19
#   row 4  -> 005046790 BEATRIZ    26166.08
#   row 5  -> 004260002 ERICA      20823.95
#   row 11 -> 004231509 THEOMAR     1021.21
# Delete from bottom to top so earlier row numbers stay valid as we go.
$rowsToDelete = @(11, 5, 4, 2)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
